$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '35.128.87'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').Value2 = '1.852.42'
$ws.Range('E3').Value = '  +1.88%  '
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '237.99'
$ws.Range('E5').Value = '  +3.38%  '
$ws.Range('E6').Value = '  +0.92%  '
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '42.13'
$ws.Range('E8').Value = '  +5.24%  '
$ws.Range('E9').Value = '  +1.38%  '
$ws.Range('E10').Value = '  +1.52%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value2 = '2.120.63'
$ws.Range('E12').Value = '  +1.93%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value2 = '1.875.29'
$ws.Range('E13').Value = '  +3.44%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '11.39'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('E15').Value = '  +0.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '4.74'
$ws.Range('E16').Value = '  +2.77%  '
$ws.Range('D17').Value2 = '35.088.06'
$ws.Range('E17').Value = '  +0.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '70.02'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '240.55'
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '12.14'
$ws.Range('E21').Value = '  +0.94%  '
$ws.Range('E22').Value = '  +2.04%  '
$ws.Range('E23').Value = '  +0.43%  '
$ws.Range('E24').Value = '  -0.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '169.90'
$ws.Range('E25').Value = '  -2.15%  '
$ws.Range('E26').Value = '  +2.37%  '
$ws.Range('E27').Value = '  +20.79%  '
$ws.Range('E28').Value = '  +1.58%  '
$ws.Range('E29').Value = '  +0.27%  '
$ws.Range('E30').Value = '  +0.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '0.0552'
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '3.98'
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('E34').Value = '  +26.38%  '
$ws.Range('E35').Value = '  +9.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '0.795'
$ws.Range('E36').Value = '  +14.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '1.29'
$ws.Range('E37').Value = '  +3.64%  '
$ws.Range('E38').Value = '  +9.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '90.19'
$ws.Range('E40').Value = '  -2.63%  '
$ws.Range('D41').Value2 = '1.345.47'
$ws.Range('E41').Value = '  +0.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '13.31'
$ws.Range('E42').Value = '  +55.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '14.95'
$ws.Range('E43').Value = '  +2.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '2.30'
$ws.Range('E44').Value = '  +1.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '2.44'
$ws.Range('E45').Value = '  +0.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '0.0556'
$ws.Range('E46').Value = '  +6.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '2.74'
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('E48').Value = '  +4.25%  '
$ws.Range('D49').Value2 = '2.028.41'
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '0.0676'
$ws.Range('E50').Value = '  +1.61%  '
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '1.01'
$ws.Range('E51').Value = '  +0.40%  '
